$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 138 (shifts existing rows 138-232 down to 139-233)
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new record
$ws.Cells.Item(138, 1).Value = 7
$ws.Cells.Item(138, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(138, 3).Value = "Ñuble"
$ws.Cells.Item(138, 4).Value = 44777
$ws.Cells.Item(138, 5).Value = 16
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100108
$ws.Cells.Item(138, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(138, 9).Value = 100108005
$ws.Cells.Item(138, 10).Value = "Piña"
$ws.Cells.Item(138, 11).Value = "Caramelo"
$ws.Cells.Item(138, 12).Value = "Segunda"
$ws.Cells.Item(138, 13).Value = 120
$ws.Cells.Item(138, 14).Value = 19000
$ws.Cells.Item(138, 15).Value = 20000
$ws.Cells.Item(138, 16).Value = 19500
$ws.Cells.Item(138, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(138, 18).Value = "Ecuador"
$ws.Cells.Item(138, 19).Value = 1393
$ws.Cells.Item(138, 20).Value = 14
